$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace measurement values for A2/D2 (was poort A1) with A2 measurement data
$ws.Range("A2").Value = 567
$ws.Range("D2").Value = 17.8

$ws.Range("A3").Value = 421
$ws.Range("D3").Value = 7.1

# Update the active selection to D4
$ws.Range("D4").Select()
